# Update CDA Logical model for ST.r2b
$wb = $excel.ActiveWorkbook

# --- Rename the "Include from CompressionAlgor" sheet to "Include #0" ---
$wsInclude = $wb.Worksheets.Item("Include from CompressionAlgor")
$wsInclude.Name = "Include #0"

# --- Metadata sheet updates ---
$ws = $wb.Worksheets.Item("Metadata")

# Bump the Version and Date property values
$ws.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"
$ws.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new "Jurisdiction" property row after "Contact" (row 10),
# pushing "Description"/"Purpose"/"Copyright"/"Immutable" down by one row.
$ws.Rows.Item(11).Insert()
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

# Give the new row the same look (font/fill/border) as the other body rows
# by copying the formatting from the row directly above it.
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false
